$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 11.2519305019305
$ws.Range("D2").Value = 9381.2101910828
$ws.Range("F2").Value = 0.860738636363636
$ws.Range("I2").Value = 49
$ws.Range("J2").Value = 264

# Row 3
$ws.Range("C3").Value = 15.2885532591415
$ws.Range("D3").Value = 4112.07349081365
$ws.Range("E3").Value = 0.0874233128834356
$ws.Range("F3").Value = 3.31141538461538
$ws.Range("G3").Value = 271
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 49
$ws.Range("J3").Value = 652

# Row 4
$ws.Range("C4").Value = 22.5581168831169
$ws.Range("D4").Value = 2148.85222381636
$ws.Range("E4").Value = 0.2868499796891
$ws.Range("F4").Value = 14.9614749475108
$ws.Range("G4").Value = 492
$ws.Range("H4").Value = 11
$ws.Range("I4").Value = 49
$ws.Range("J4").Value = 1189

# Row 5
$ws.Range("C5").Value = 32.1375530410184
$ws.Range("D5").Value = 1285.99128540305
$ws.Range("E5").Value = 14.1675338453935
$ws.Range("F5").Value = 27.8125796114794
$ws.Range("G5").Value = 274
$ws.Range("H5").Value = 16
$ws.Range("I5").Value = 49
$ws.Range("J5").Value = 733

# Row 6
$ws.Range("C6").Value = 14.183859223301
$ws.Range("D6").Value = 14140.9217877095
$ws.Range("F6").Value = 0.346694553920201
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 45
$ws.Range("J6").Value = 211

# Row 7
$ws.Range("C7").Value = 19.25725
$ws.Range("D7").Value = 4959.55497382199
$ws.Range("E7").Value = 0.0405339805825243
$ws.Range("F7").Value = 5.56043164021548
$ws.Range("G7").Value = 133
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 45
$ws.Range("J7").Value = 515

# Row 8
$ws.Range("C8").Value = 25.6649045020464
$ws.Range("D8").Value = 2319.40104166667
$ws.Range("E8").Value = 0.585194909621996
$ws.Range("F8").Value = 23.2136923109855
$ws.Range("G8").Value = 181
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 45
$ws.Range("J8").Value = 757

# Row 9
$ws.Range("C9").Value = 36.4011627906977
$ws.Range("D9").Value = 1297.90273556231
$ws.Range("E9").Value = 14.7327797866167
$ws.Range("F9").Value = 45.7196431269731
$ws.Range("G9").Value = 64
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 45
$ws.Range("J9").Value = 393

# Row 14
$ws.Range("F14").Value = 0.567888133728462
$ws.Range("H14").Value = 0

# Row 15
$ws.Range("F15").Value = 1.20445965009235
$ws.Range("H15").Value = 1

# Row 16
$ws.Range("F16").Value = 4.75284541709586
$ws.Range("H16").Value = 2

# Row 17
$ws.Range("E17").Value = 4.2484037050403
$ws.Range("F17").Value = 25.2876374348006
$ws.Range("H17").Value = 1
